$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the duplicate "Polynesia" row (erroneous standalone region entry;
# "French Polynesia" already carries the PYF code on the next row).
$ws.Rows.Item(124).Delete()

# Keep the active selection in sync with the new bottom of the shifted table.
$ws.Range("B458").Select()
